$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(6, 20, 20, 20, 1, 2.5, 1, 0.5, 0, 15, 10, 0, 0, 0, 0.01, 0.01, -0.03954538461277818, 26.62025762178805, 16.31018336053107, 0, "eyescrunching"),
    @(6, 20, 20, 20, 1, 2.5, 1, 0.5, 0, 15, 10, 0, 0, 0, 0.01, 0.01, -0.02081399515663128, 23.06462261758369, 21.83555164121976, 1, "jaw"),
    @(6, 20, 20, 20, 1, 2.5, 1, 0.5, 0, 15, 10, 0, 0, 0, 0.01, 0.01, -0.03782867504932325, 16.43695554907289, -14.5564680488108, 2, "raisingeyebrows"),
    @(6, 20, 20, 20, 1, 2.5, 1, 0.5, 0, 15, 10, 0, 0, 0, 0.001, 0.001, -3.04077830408176, 9.925443034519448, -4.67507000742246, 3, "movehat"),
    @(6, 20, 20, 20, 1, 2.5, 1, 0.5, 0, 15, 10, 0, 0, 0, 0.01, 0.1, -0.03688754513089874, 1.56427737596919, -8.457205793158687, 4, "movehead")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($c = 1; $c -le $values.Length; $c++) {
        $ws.Cells.Item($row, $c).Value = $values[$c - 1]
    }
}
